# Auto-committed on 2021/12/30 週四
# Insert a new "ServerIp" row into the DBD field list (TxPrinter table),
# renumber the SEQ column with formulas, update the "Printer" row's
# Chinese comment, and make DBD the active/selected sheet again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws2 = $wb.Worksheets.Item("DBS")

# --- insert a new row 11 (pushes old rows 11-16 down to 12-17) ---------
$ws.Rows("11:11").Insert()

# Rebuild per-cell styles on the new row 11 by pulling formats from the
# two donor rows that already carry the exact look we need: row 9
# (StanIp) supplies A/C/D/E/F/G, row 10 (FileCode) supplies B.
$ws.Range("A9:G9").Copy()
$ws.Range("A11:G11").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- new row 11 content: ServerIp field ---------------------------------
$ws.Range("B11").Value = "ServerIp"
$ws.Range("C11").Value = "印表機伺服器IP"
$ws.Range("D11").Value = "varchar2"
$ws.Range("E11").Value = "15"

# --- row 12 (shifted from the old row 11) keeps the Printer field but --
# --- gets a new Chinese comment -----------------------------------------
$ws.Range("C12").Value = "預設印表機"

# --- SEQ column: turn A10:A16 into the shared "+1" formula chain -------
$ws.Range("A11:A16").Formula = "=A10+1"
$ws.Range("A10").Formula = "=A9+1"

# --- make DBD the active/selected sheet again, with C12 selected -------
$ws.Activate()
$ws.Range("C12").Select()
